$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2375.6765
$ws.Range("I132").Value = 2459.1538
$ws.Range("J132").Value = 2104.375
$ws.Range("K132").Value = 7377.4614
$ws.Range("L132").Value = 6313.125
$ws.Range("M132").Value = -4847.4614
$ws.Range("N132").Value = -11373.125
$ws.Range("H135").Value = 14707323
$ws.Range("I135").Value = 980.9
$ws.Range("K135").Value = 8828.1
$ws.Range("M135").Value = -6293.1
$ws.Range("H137").Value = 2315.8108
$ws.Range("I137").Value = 2439.6924
$ws.Range("K137").Value = 7319.0772
$ws.Range("M137").Value = -4769.0772

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11369378
$ws.Range("I32").Value = 14708045
$ws.Range("J32").Value = 17912.3
$ws.Range("K32").Value = 14708045
$ws.Range("L32").Value = 17912.3
$ws.Range("M32").Value = -14707758
$ws.Range("N32").Value = -18486.3
$ws.Range("H45").Value = 1652.6666
$ws.Range("I45").Value = 1237.3334
$ws.Range("K45").Value = 1237.3334
$ws.Range("M45").Value = -860.3334
$ws.Range("H61").Value = 29414712
$ws.Range("I61").Value = 41667864
$ws.Range("K61").Value = 41667864
$ws.Range("M61").Value = -41667652
$ws.Range("H110").Value = 11337.541
$ws.Range("I110").Value = 11918.303
$ws.Range("J110").Value = 6546.25
$ws.Range("K110").Value = 11918.303
$ws.Range("L110").Value = 6546.25
$ws.Range("M110").Value = -9873.303
$ws.Range("N110").Value = -10636.25
$ws.Range("H136").Value = 29414712
$ws.Range("I136").Value = 41667864
$ws.Range("K136").Value = 125003592
$ws.Range("M136").Value = -125001042
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2008.8636
$ws.Range("I20").Value = 1982.3125
$ws.Range("K20").Value = 1982.3125
$ws.Range("M20").Value = -1735.3125
$ws.Range("H105").Value = 8673.105
$ws.Range("I105").Value = 13888.556
$ws.Range("K105").Value = 13888.556
$ws.Range("M105").Value = -12141.556
$ws.Range("H134").Value = 3940.6897
$ws.Range("I134").Value = 3862.2593
$ws.Range("K134").Value = 11586.7779
$ws.Range("M134").Value = -9051.777900000001
$ws.Range("H141").Value = 89000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 89000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 89000
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -99360

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5985.2104
$ws.Range("I22").Value = 7054
$ws.Range("J22").Value = 1977.25
$ws.Range("K22").Value = 7054
$ws.Range("L22").Value = 1977.25
$ws.Range("M22").Value = -6704
$ws.Range("N22").Value = -2677.25
$ws.Range("H31").Value = 22731930
$ws.Range("I31").Value = 3137.0557
$ws.Range("J31").Value = 125011496
$ws.Range("K31").Value = 3137.0557
$ws.Range("L31").Value = 125011496
$ws.Range("M31").Value = -2842.0557
$ws.Range("N31").Value = -125012086
$ws.Range("H34").Value = 22731930
$ws.Range("I34").Value = 3137.0557
$ws.Range("J34").Value = 125011496
$ws.Range("K34").Value = 3137.0557
$ws.Range("L34").Value = 125011496
$ws.Range("M34").Value = -2935.0557
$ws.Range("N34").Value = -125011900
$ws.Range("H58").Value = 3559.6667
$ws.Range("I58").Value = 2939.5
$ws.Range("J58").Value = 4800
$ws.Range("K58").Value = 2939.5
$ws.Range("L58").Value = 4800
$ws.Range("M58").Value = -2736.5
$ws.Range("N58").Value = -5206
$ws.Range("H99").Value = 12145.605
$ws.Range("I99").Value = 10061.9375
$ws.Range("J99").Value = 13661
$ws.Range("K99").Value = 10061.9375
$ws.Range("L99").Value = 13661
$ws.Range("M99").Value = -8563.9375
$ws.Range("N99").Value = -16657
$ws.Range("H126").Value = 12145.605
$ws.Range("I126").Value = 10061.9375
$ws.Range("J126").Value = 13661
$ws.Range("K126").Value = 30185.8125
$ws.Range("L126").Value = 40983
$ws.Range("M126").Value = -27715.8125
$ws.Range("N126").Value = -45923
$ws.Range("H136").Value = 3559.6667
$ws.Range("I136").Value = 2939.5
$ws.Range("J136").Value = 4800
$ws.Range("K136").Value = 8818.5
$ws.Range("L136").Value = 14400
$ws.Range("M136").Value = -6268.5
$ws.Range("N136").Value = -19500
$ws.Range("H141").Value = 301214.16
$ws.Range("J141").Value = 376586.12
$ws.Range("L141").Value = 376586.12
$ws.Range("N141").Value = -386946.12

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 161663.67
$ws.Range("J37").Value = 161663.67
$ws.Range("L37").Value = 484991.01
$ws.Range("N37").Value = -485215.01
$ws.Range("H86").Value = 677.4545000000001
$ws.Range("I86").Value = 478
$ws.Range("K86").Value = 1434
$ws.Range("M86").Value = -248
$ws.Range("H89").Value = 677.4545000000001
$ws.Range("I89").Value = 478
$ws.Range("K89").Value = 4302
$ws.Range("M89").Value = 1626
$ws.Range("H113").Value = 3035.3845
$ws.Range("I113").Value = 2899.6667
$ws.Range("J113").Value = 3076.1
$ws.Range("K113").Value = 8699.000100000001
$ws.Range("L113").Value = 9228.299999999999
$ws.Range("M113").Value = -6529.000100000001
$ws.Range("N113").Value = -13568.3
$ws.Range("H140").Value = 1650.5
$ws.Range("I140").Value = 1650.5
$ws.Range("K140").Value = 4951.5
$ws.Range("M140").Value = 228.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1132.7059
$ws.Range("I97").Value = 393.25
$ws.Range("J97").Value = 2907.4
$ws.Range("K97").Value = 393.25
$ws.Range("L97").Value = 2907.4
$ws.Range("M97").Value = 102.75
$ws.Range("N97").Value = -3899.4
$ws.Range("H132").Value = 2387.5557
$ws.Range("I132").Value = 2367.6086
$ws.Range("J132").Value = 2502.25
$ws.Range("K132").Value = 7102.825800000001
$ws.Range("L132").Value = 7506.75
$ws.Range("M132").Value = -4572.825800000001
$ws.Range("N132").Value = -12566.75

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2782.5
$ws.Range("I68").Value = 2536.5625
$ws.Range("J68").Value = 4750
$ws.Range("K68").Value = 2536.5625
$ws.Range("L68").Value = 4750
$ws.Range("M68").Value = -1787.5625
$ws.Range("N68").Value = -6248
$ws.Range("H71").Value = 2782.5
$ws.Range("I71").Value = 2536.5625
$ws.Range("J71").Value = 4750
$ws.Range("K71").Value = 12682.8125
$ws.Range("L71").Value = 23750
$ws.Range("M71").Value = -8938.8125
$ws.Range("N71").Value = -31238
$ws.Range("H82").Value = 2370.8215
$ws.Range("I82").Value = 1405
$ws.Range("K82").Value = 1405
$ws.Range("M82").Value = -1044
$ws.Range("H85").Value = 2370.8215
$ws.Range("I85").Value = 1405
$ws.Range("K85").Value = 1405
$ws.Range("M85").Value = -157
$ws.Range("H122").Value = 3706.625
$ws.Range("I122").Value = 3164.7144
$ws.Range("J122").Value = 7500
$ws.Range("K122").Value = 9494.143199999999
$ws.Range("L122").Value = 22500
$ws.Range("M122").Value = -7044.143199999999
$ws.Range("N122").Value = -27400
$ws.Range("H131").Value = 88992
$ws.Range("J131").Value = 88992
$ws.Range("L131").Value = 88992
$ws.Range("N131").Value = -99072
$ws.Range("H132").Value = 90912430
$ws.Range("I132").Value = 3045.5334
$ws.Range("J132").Value = 285718270
$ws.Range("K132").Value = 9136.600199999999
$ws.Range("L132").Value = 857154810
$ws.Range("M132").Value = -6606.600199999999
$ws.Range("N132").Value = -857159870
$ws.Range("H136").Value = 2234.0312
$ws.Range("I136").Value = 1696.4286
$ws.Range("K136").Value = 5089.2858
$ws.Range("M136").Value = -2539.2858

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 45501500
$ws.Range("I122").Value = 50051200
$ws.Range("K122").Value = 150153600
$ws.Range("M122").Value = -150151150
$ws.Range("H126").Value = 3024.1177
$ws.Range("I126").Value = 3024.1177
$ws.Range("K126").Value = 9072.3531
$ws.Range("M126").Value = -6602.3531
$ws.Range("H132").Value = 5004.512
$ws.Range("I132").Value = 4911.6665
$ws.Range("J132").Value = 5673
$ws.Range("K132").Value = 14734.9995
$ws.Range("L132").Value = 17019
$ws.Range("M132").Value = -12204.9995
$ws.Range("N132").Value = -22079
$ws.Range("H136").Value = 1714.7241
$ws.Range("I136").Value = 1215.76
$ws.Range("J136").Value = 4833.25
$ws.Range("K136").Value = 3647.28
$ws.Range("L136").Value = 14499.75
$ws.Range("M136").Value = -1097.28
$ws.Range("N136").Value = -19599.75
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
